$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of rows 215-227 (columns A:C) while keeping formatting
$ws.Range("A215:C227").ClearContents()

# Update the active selection to E215 as shown in the diff
$ws.Range("E215").Select()
